$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.521.92"
$ws.Range("E2").Value = "  +2.52%  "

$ws.Range("D3").Value = "2.076.22"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'235.11"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  +3.99%  "

$ws.Range("D7").Value = "'58.27"
$ws.Range("E7").Value = "  +6.45%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +4.36%  "

$ws.Range("D10").Value = "'59.13"
$ws.Range("E10").Value = "  +3.30%  "

$ws.Range("E11").Value = "  +2.49%  "

$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").Value = "2.383.59"
$ws.Range("E13").Value = "  +3.57%  "

$ws.Range("D14").Value = "'14.49"
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "'21.07"
$ws.Range("E15").Value = "  +4.94%  "

$ws.Range("E16").Value = "  +3.67%  "

$ws.Range("D17").Value = "'5.19"
$ws.Range("E17").Value = "  +2.18%  "

$ws.Range("D18").Value = "2.087.74"
$ws.Range("E18").Value = "  +4.01%  "

$ws.Range("D19").Value = "37.721.88"
$ws.Range("E19").Value = "  +3.40%  "

$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  +16.79%  "

$ws.Range("D21").Value = "'69.55"
$ws.Range("E21").Value = "  +2.85%  "

$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("D23").Value = "'226.26"
$ws.Range("E23").Value = "  +2.38%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "'2.51"
$ws.Range("E25").Value = "  +4.81%  "

$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").Value = "'167.13"
$ws.Range("E27").Value = "  +2.57%  "

$ws.Range("E28").Value = "  +7.69%  "

$ws.Range("E29").Value = "  +4.67%  "

$ws.Range("E30").Value = "  +2.46%  "

$ws.Range("D31").Value = "'19.25"
$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("E32").Value = "  +1.89%  "

$ws.Range("D33").Value = "'4.55"
$ws.Range("E33").Value = "  +4.11%  "

$ws.Range("E34").Value = "  +4.07%  "

$ws.Range("D35").Value = "'2.58"
$ws.Range("E35").Value = "  +7.53%  "

$ws.Range("D36").Value = "'4.59"
$ws.Range("E36").Value = "  +8.46%  "

$ws.Range("D37").Value = "'6.00"
$ws.Range("E37").Value = "  +4.84%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").Value = "'4.59"
$ws.Range("E41").Value = "  +19.03%  "

$ws.Range("D42").Value = "'2.95"
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("D43").Value = "'0.0952"
$ws.Range("E43").Value = "  +3.08%  "

$ws.Range("D44").Value = "1.469.48"
$ws.Range("E44").Value = "  +0.92%  "

$ws.Range("E45").Value = "  +7.09%  "

$ws.Range("D46").Value = "'95.62"
$ws.Range("E46").Value = "  +6.61%  "

$ws.Range("E47").Value = "  +4.36%  "

$ws.Range("D48").Value = "'15.88"
$ws.Range("E48").Value = "  +3.86%  "

$ws.Range("E49").Value = "  +3.94%  "

$ws.Range("D50").Value = "'7.20"
$ws.Range("E50").Value = "  +5.37%  "

$ws.Range("E51").Value = "  +2.13%  "
